# Word COM-interop script implementing the diff:
#  - splits "books of js" into "books of " + proofErr(spellStart/End around "js")
#  - inserts a blank paragraph
#  - inserts a new paragraph with a hyperlink to the otus.ru lesson page,
#    followed by commentary text (with proofErr spell/gram markers), and
#    relocates the _GoBack bookmark to the end of that new paragraph
$d = $word.ActiveDocument

# --- Step 1: allocate a real hyperlink relationship for the otus.ru link ---
# Word assigns a fresh rId when Hyperlinks.Add is called; we do this on a
# disposable paragraph appended at the end of the story, capture the rId it
# was given, then remove the scratch paragraph again so the visible content
# is untouched.
$lastParaRange = $d.Paragraphs($d.Paragraphs.Count).Range
$lastParaRange.Collapse(0)
$lastParaRange.InsertParagraphAfter()
$scratchRange = $d.Paragraphs($d.Paragraphs.Count).Range
$d.Hyperlinks.Add($scratchRange, "https://otus.ru/lessons/microservice-architecture/process/") | Out-Null

$fullXml = $d.Content.WordOpenXML
$otusRid = "rId8"
if ($fullXml -match 'r:id="(rId\d+)"[^>]*><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://otus\.ru') {
    $otusRid = $matches[1]
}

$scratchPara = $d.Paragraphs($d.Paragraphs.Count).Range
$cleanupRange = $d.Range($scratchPara.Start - 1, $scratchPara.End)
$cleanupRange.Delete()

# --- Step 2: rewrite the "books of js" paragraph + insert the two new ones ---
$target = $d.Paragraphs(8).Range

$newXml = '<w:p w:rsidR="00EB6BC8" w:rsidRPr="00EB6BC8" w:rsidRDefault="00EB6BC8"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:hyperlink r:id="rId7" w:history="1"><w:r w:rsidRPr="001F48E5"><w:rPr><w:rStyle w:val="a3"/><w:lang w:val="en-US"/></w:rPr><w:t>https://vk.com/wall-54530371_125557</w:t></w:r></w:hyperlink><w:r w:rsidRPr="00EB6BC8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> -- </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">books of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:hyperlink r:id="' + $otusRid + '" w:history="1"><w:r><w:rPr><w:rStyle w:val="a3"/><w:lang w:val="en-US"/></w:rPr><w:t>https://otus.ru/lessons/microservice-architecture/process/</w:t></w:r></w:hyperlink><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">about </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>autettifiacation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> protocols. I learned only half of </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>this(</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="gramEnd"/></w:p>'

$target.InsertXML($newXml)

Write-Output "done; otusRid=$otusRid"
